$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-23 00:16:11"
$overview.Range("G3").Value = "2016-08-23 00:16:11"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H2").Value = "2016-08-23 00:16:01"
$zhcn.Range("H3").Value = "2016-08-23 00:16:01"
$zhcn.Range("K2").Value = "2016-08-23 00:16:28"
$zhcn.Range("K3").Value = "2016-08-23 00:16:28"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "mt"
$dede.Range("E3").Value = "mt"
$dede.Range("H2").Value = "2016-08-23 00:16:11"
$dede.Range("H3").Value = "2016-08-23 00:16:11"
$dede.Range("K2").Value = "2016-08-23 00:16:35"
$dede.Range("K3").Value = "2016-08-23 00:16:35"
